$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended to the bottom of the chat-log sheet (rows 47 & 48),
# matching the "工众号" / "并音工众号" additions from the diff.
$ws.Range("A47").Value = "工众号"
$ws.Range("B47").Value = 1
$ws.Range("A48").Value = "并音工众号"
$ws.Range("B48").Value = 1

# Match the existing style used for the header row / body text cells
# (vertical-center alignment -> style index 1 in styles.xml).
$ws.Range("A47:A48").VerticalAlignment = -4108

# Move the view/selection down to the newly-added rows, like the author
# scrolling down after typing the new entries.
$ws.Application.Goto($ws.Range("A47"), $false)
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A47").Select()
